$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10 is written first so its label ("CONFIDENCE.T") lands at shared-string
# index 9 and row 9's label ("CONFIDENCE.NORM") lands at index 10 - matching
# the insertion order used by the original author.
$ws.Range("A10").Value = "CONFIDENCE.T"
$ws.Range("B10").Formula = "=_xlfn.CONFIDENCE.NORM(C10,D10,E10)"
$ws.Range("C10").Formula = "=2/15"
$ws.Range("D10").Value = 6.6
$ws.Range("E10").Value = 44
$ws.Range("B10").NumberFormat = "0.0000000000"

$ws.Range("A9").Value = "CONFIDENCE.NORM"
$ws.Range("B9").Formula = "=_xlfn.CONFIDENCE.NORM(C9,D9,E9)"
$ws.Range("C9").Formula = "=2/15"
$ws.Range("D9").Value = 6.6
$ws.Range("E9").Value = 44
$ws.Range("B9").NumberFormat = "0.0000000000"

# Column A widened to fit the new longer labels.
$ws.Columns.Item(1).ColumnWidth = 18.5

# Selection moves on to the next empty row, like Excel leaves it after data entry.
$ws.Range("A11").Select()
